$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the "Program Definitions" and "Metadata" sheets
$wb.Worksheets.Item("Program Definitions").Delete()
$wb.Worksheets.Item("Metadata").Delete()

# Update population-name lookups in "Parameters" to reference the
# abbreviation column (A) instead of the full-name column (B)
$paramRows = @(2, 3, 6, 7, 10, 11, 14, 15, 18, 19)
$wsParams = $wb.Worksheets.Item("Parameters")
foreach ($r in $paramRows) {
    $srcRow = 2 + (($r - 2) % 2)
    $wsParams.Range("A$r").Formula = "='Population Definitions'!A$srcRow"
}

# Same update for "State Variables"
$stateRows = @(2, 3, 6, 7, 10, 11)
$wsState = $wb.Worksheets.Item("State Variables")
foreach ($r in $stateRows) {
    $srcRow = 2 + (($r - 2) % 2)
    $wsState.Range("A$r").Formula = "='Population Definitions'!A$srcRow"
}

# Restore cursor/selection positions on each remaining sheet
$wsPop = $wb.Worksheets.Item("Population Definitions")
$wsPop.Range("B2").Select()

$wsParams.Range("A2").Select()

$wsState.Activate()
$wsState.Range("X13").Select()

Write-Host "Done"
